# The source workbook contained an accidental duplicate data row: row 131
# ("Syria") duplicated the data already present in row 127. This edit
# removes that duplicate row, which shifts every row below it (the last
# five "extreme alarming" countries -- Eritrea, Libya, Maldives, Qatar --
# and the now-departed second Syria line) up by one, and shrinks the
# sheet's used range from A1:K135 to A1:K134.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2022 GHI Ranking - Tableau 1")

# Remove the duplicated "Syria" row.
$ws.Rows("131:131").Delete()

# Row deletion shrinks the used range by one row (now A1:K134), but the
# sheet's AutoFilter and its backing _FilterDatabase defined name are
# "sticky" and keep pointing at the old A1:E135 extent, so bring them
# back in sync with the new data extent explicitly.
$ws.AutoFilterMode = $false
$ws.Range("A1:E134").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='2022 GHI Ranking - Tableau 1'!`$A`$1:`$E`$134"
    }
}

# Leave the view where the user ended up after performing the cleanup:
# scrolled down to the tail of the table with the next duplicate-looking
# row selected.
$ws.Range("B132").Select()
